# First phase of implementing EAF (Evicted Address Filter) Policy
#
# Updates the "Results-Without Warmup-WorseCase" sheet:
#   - Refreshed Total Misses (E column, rows 6-16) and CPI (E column, rows 38-48)
#     numbers for the "Contestant Policy" column, now that DRRIP defines/knobs
#     have been wired in.
#   - New FC-policy numbers land in F11 / F43 (previously empty).
#   - New column I holds a short dump of the DRRIP #define block used to
#     produce the contestant-policy run.
#   - Re-point the view/selection at the area that was being edited.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Total Misses (Worse Case) : Contestant Policy column (E), rows 6-16 ----
$ws.Range("E6").Value  = 11140507
$ws.Range("E7").Value  = 173404
$ws.Range("E8").Value  = 469623
$ws.Range("E9").Value  = 940602
$ws.Range("E10").Value = 931812
$ws.Range("E11").Value = 105887
$ws.Range("F11").Value = 75950
$ws.Range("E12").Value = 1432373
$ws.Range("E13").Value = 2217199
$ws.Range("E14").Value = 10271704
$ws.Range("E15").Value = 1424331
$ws.Range("E16").Value = 10448328

# ---- CPI (Worse Case) : Contestant Policy column (E), rows 38-48 ----
$ws.Range("E38").Value = 4.43992
$ws.Range("E39").Value = 0.289877
$ws.Range("E40").Value = 0.525547
$ws.Range("E41").Value = 0.869599
$ws.Range("E42").Value = 0.659096
$ws.Range("E43").Value = 0.294995
$ws.Range("F43").Value = 0.284656
$ws.Range("E44").Value = 0.902738
$ws.Range("E45").Value = 1.48793
$ws.Range("E46").Value = 4.17935
$ws.Range("E47").Value = 0.8996
$ws.Range("E48").Value = 4.22662

# ---- New column I : DRRIP defines notes, rows 6-10 ----
$ws.Range("I6").Value  = "//DRRIP Defines"
$ws.Range("I7").Value  = "#define NumLeaderSets   64"
$ws.Range("I8").Value  = "#define RRIP_MAX        3"
$ws.Range("I9").Value  = "#define PSEL_MAX        15"
$ws.Range("I10").Value = "#define BIOMODAL_PROBABILITY    31   //[1 means 0.1%/10 means 1%] of all times"

# ---- View / selection: focus back on the area just edited ----
$ws.Range("F43").Select()
